# Auto-generated edit script for "Add data for 2022-07-15"
# Updates violent crime counts (mostly the 2022 / column I totals) across
# the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 9).Value = 3684   # I2: 3665 -> 3684
$ws.Cells.Item(3, 7).Value = 8320   # G3: 8319 -> 8320
$ws.Cells.Item(3, 9).Value = 3802   # I3: 3791 -> 3802
$ws.Cells.Item(4, 6).Value = 1860   # F4: 1859 -> 1860
$ws.Cells.Item(4, 9).Value = 893   # I4: 892 -> 893
$ws.Cells.Item(5, 9).Value = 350   # I5: 349 -> 350
$ws.Cells.Item(6, 9).Value = 4275   # I6: 4255 -> 4275
$ws.Cells.Item(7, 6).Value = 24049   # F7: 24048 -> 24049
$ws.Cells.Item(7, 7).Value = 24659   # G7: 24658 -> 24659
$ws.Cells.Item(7, 9).Value = 13004   # I7: 12952 -> 13004

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item(14)
$ws.Cells.Item(3, 9).Value = 128   # I3: 127 -> 128
$ws.Cells.Item(6, 9).Value = 116   # I6: 115 -> 116
$ws.Cells.Item(7, 9).Value = 418   # I7: 416 -> 418

# Sheet 15: Woodlawn
$ws = $wb.Worksheets.Item(15)
$ws.Cells.Item(2, 9).Value = 64   # I2: 63 -> 64
$ws.Cells.Item(3, 9).Value = 87   # I3: 86 -> 87
$ws.Cells.Item(4, 9).Value = 19   # I4: 18 -> 19
$ws.Cells.Item(6, 9).Value = 68   # I6: 67 -> 68
$ws.Cells.Item(7, 9).Value = 244   # I7: 240 -> 244

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item(16)
$ws.Cells.Item(2, 9).Value = 121   # I2: 120 -> 121
$ws.Cells.Item(6, 9).Value = 166   # I6: 167 -> 166

# Sheet 18: South Deering
$ws = $wb.Worksheets.Item(18)
$ws.Cells.Item(2, 9).Value = 42   # I2: 41 -> 42
$ws.Cells.Item(3, 9).Value = 37   # I3: 35 -> 37
$ws.Cells.Item(7, 9).Value = 114   # I7: 111 -> 114

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 9).Value = 114   # I2: 113 -> 114
$ws.Cells.Item(6, 9).Value = 89   # I6: 88 -> 89
$ws.Cells.Item(7, 9).Value = 421   # I7: 419 -> 421
$ws.Cells.Item(8, 9).Value = 790   # I8: 787 -> 790
$ws.Cells.Item(10, 9).Value = 89   # I10: 88 -> 89
$ws.Cells.Item(11, 9).Value = 204   # I11: 203 -> 204
$ws.Cells.Item(15, 9).Value = 155   # I15: 154 -> 155
$ws.Cells.Item(18, 9).Value = 88   # I18: 87 -> 88
$ws.Cells.Item(19, 9).Value = 347   # I19: 345 -> 347
$ws.Cells.Item(21, 9).Value = 71   # I21: 70 -> 71
$ws.Cells.Item(29, 9).Value = 840   # I29: 841 -> 840
$ws.Cells.Item(33, 9).Value = 584   # I33: 581 -> 584
$ws.Cells.Item(36, 9).Value = 179   # I36: 178 -> 179
$ws.Cells.Item(37, 9).Value = 418   # I37: 416 -> 418
$ws.Cells.Item(41, 9).Value = 58   # I41: 57 -> 58
$ws.Cells.Item(42, 9).Value = 454   # I42: 452 -> 454
$ws.Cells.Item(43, 9).Value = 114   # I43: 112 -> 114
$ws.Cells.Item(47, 9).Value = 91   # I47: 90 -> 91
$ws.Cells.Item(48, 9).Value = 170   # I48: 167 -> 170
$ws.Cells.Item(51, 9).Value = 125   # I51: 124 -> 125
$ws.Cells.Item(52, 9).Value = 283   # I52: 282 -> 283
$ws.Cells.Item(54, 9).Value = 294   # I54: 291 -> 294
$ws.Cells.Item(63, 7).Value = 200   # G63: 199 -> 200
$ws.Cells.Item(63, 9).Value = 51   # I63: 54 -> 51
$ws.Cells.Item(64, 9).Value = 115   # I64: 117 -> 115
$ws.Cells.Item(70, 9).Value = 24   # I70: 23 -> 24
$ws.Cells.Item(78, 9).Value = 185   # I78: 184 -> 185
$ws.Cells.Item(79, 9).Value = 349   # I79: 348 -> 349
$ws.Cells.Item(83, 9).Value = 264   # I83: 262 -> 264
$ws.Cells.Item(84, 9).Value = 114   # I84: 111 -> 114
$ws.Cells.Item(85, 9).Value = 590   # I85: 587 -> 590
$ws.Cells.Item(86, 9).Value = 78   # I86: 77 -> 78
$ws.Cells.Item(90, 9).Value = 164   # I90: 162 -> 164
$ws.Cells.Item(91, 9).Value = 160   # I91: 158 -> 160
$ws.Cells.Item(93, 9).Value = 74   # I93: 73 -> 74
$ws.Cells.Item(94, 9).Value = 120   # I94: 119 -> 120
$ws.Cells.Item(95, 9).Value = 206   # I95: 205 -> 206
$ws.Cells.Item(97, 6).Value = 156   # F97: 155 -> 156
$ws.Cells.Item(97, 9).Value = 93   # I97: 92 -> 93
$ws.Cells.Item(99, 9).Value = 244   # I99: 240 -> 244
$ws.Cells.Item(101, 6).Value = 24049   # F101: 24048 -> 24049
$ws.Cells.Item(101, 7).Value = 24659   # G101: 24658 -> 24659
$ws.Cells.Item(101, 9).Value = 13004   # I101: 12952 -> 13004

# Sheet 20: South Chicago
$ws = $wb.Worksheets.Item(20)
$ws.Cells.Item(2, 9).Value = 92   # I2: 91 -> 92
$ws.Cells.Item(6, 9).Value = 48   # I6: 47 -> 48
$ws.Cells.Item(7, 9).Value = 264   # I7: 262 -> 264

# Sheet 21: West Pullman
$ws = $wb.Worksheets.Item(21)
$ws.Cells.Item(2, 9).Value = 73   # I2: 72 -> 73
$ws.Cells.Item(7, 9).Value = 206   # I7: 205 -> 206

# Sheet 22: Garfield Park
$ws = $wb.Worksheets.Item(22)
$ws.Cells.Item(2, 9).Value = 139   # I2: 138 -> 139
$ws.Cells.Item(6, 9).Value = 187   # I6: 185 -> 187
$ws.Cells.Item(7, 9).Value = 584   # I7: 581 -> 584

# Sheet 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Cells.Item(2, 9).Value = 67   # I2: 65 -> 67
$ws.Cells.Item(3, 9).Value = 57   # I3: 56 -> 57
$ws.Cells.Item(7, 9).Value = 294   # I7: 291 -> 294

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Cells.Item(3, 9).Value = 286   # I3: 287 -> 286
$ws.Cells.Item(4, 9).Value = 38   # I4: 40 -> 38
$ws.Cells.Item(5, 9).Value = 33   # I5: 32 -> 33
$ws.Cells.Item(6, 9).Value = 232   # I6: 231 -> 232
$ws.Cells.Item(7, 9).Value = 840   # I7: 841 -> 840

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item(26)
$ws.Cells.Item(3, 9).Value = 100   # I3: 99 -> 100
$ws.Cells.Item(4, 9).Value = 13   # I4: 12 -> 13
$ws.Cells.Item(7, 9).Value = 347   # I7: 345 -> 347

# Sheet 28: Lake View
$ws = $wb.Worksheets.Item(28)
$ws.Cells.Item(4, 9).Value = 17   # I4: 16 -> 17
$ws.Cells.Item(6, 9).Value = 97   # I6: 95 -> 97
$ws.Cells.Item(7, 9).Value = 170   # I7: 167 -> 170

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 9).Value = 155   # I2: 153 -> 155
$ws.Cells.Item(3, 9).Value = 236   # I3: 235 -> 236
$ws.Cells.Item(7, 9).Value = 590   # I7: 587 -> 590

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item(30)
$ws.Cells.Item(2, 9).Value = 40   # I2: 39 -> 40
$ws.Cells.Item(7, 9).Value = 89   # I7: 88 -> 89

# Sheet 31: Hermosa
$ws = $wb.Worksheets.Item(31)
$ws.Cells.Item(4, 9).Value = 6   # I4: 5 -> 6
$ws.Cells.Item(7, 9).Value = 58   # I7: 57 -> 58

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item(32)
$ws.Cells.Item(2, 9).Value = 124   # I2: 123 -> 124
$ws.Cells.Item(6, 9).Value = 123   # I6: 122 -> 123
$ws.Cells.Item(7, 9).Value = 454   # I7: 452 -> 454

# Sheet 34: Avondale
$ws = $wb.Worksheets.Item(34)
$ws.Cells.Item(3, 9).Value = 17   # I3: 16 -> 17
$ws.Cells.Item(7, 9).Value = 89   # I7: 88 -> 89

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item(35)
$ws.Cells.Item(6, 9).Value = 75   # I6: 74 -> 75
$ws.Cells.Item(7, 9).Value = 185   # I7: 184 -> 185

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item(40)
$ws.Cells.Item(4, 9).Value = 4   # I4: 3 -> 4
$ws.Cells.Item(6, 9).Value = 51   # I6: 50 -> 51
$ws.Cells.Item(7, 9).Value = 160   # I7: 158 -> 160

# Sheet 41: Chinatown
$ws = $wb.Worksheets.Item(41)
$ws.Cells.Item(3, 9).Value = 9   # I3: 8 -> 9
$ws.Cells.Item(7, 9).Value = 71   # I7: 70 -> 71

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Cells.Item(6, 9).Value = 104   # I6: 103 -> 104
$ws.Cells.Item(7, 9).Value = 349   # I7: 348 -> 349

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item(43)
$ws.Cells.Item(6, 9).Value = 41   # I6: 43 -> 41
$ws.Cells.Item(7, 9).Value = 115   # I7: 117 -> 115

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item(45)
$ws.Cells.Item(2, 9).Value = 25   # I2: 24 -> 25
$ws.Cells.Item(7, 9).Value = 88   # I7: 87 -> 88

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item(47)
$ws.Cells.Item(6, 9).Value = 54   # I6: 53 -> 54
$ws.Cells.Item(7, 9).Value = 179   # I7: 178 -> 179

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item(48)
$ws.Cells.Item(6, 9).Value = 28   # I6: 27 -> 28
$ws.Cells.Item(7, 9).Value = 74   # I7: 73 -> 74

# Sheet 5: Little Village
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 9).Value = 29   # I4: 28 -> 29
$ws.Cells.Item(7, 9).Value = 283   # I7: 282 -> 283

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item(51)
$ws.Cells.Item(6, 9).Value = 67   # I6: 66 -> 67
$ws.Cells.Item(7, 9).Value = 120   # I7: 119 -> 120

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item(53)
$ws.Cells.Item(2, 9).Value = 15   # I2: 14 -> 15
$ws.Cells.Item(7, 9).Value = 91   # I7: 90 -> 91

# Sheet 54: Brighton Park
$ws = $wb.Worksheets.Item(54)
$ws.Cells.Item(2, 9).Value = 49   # I2: 48 -> 49
$ws.Cells.Item(7, 9).Value = 155   # I7: 154 -> 155

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(6, 9).Value = 47   # I6: 46 -> 47
$ws.Cells.Item(7, 9).Value = 204   # I7: 203 -> 204

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item(64)
$ws.Cells.Item(6, 9).Value = 23   # I6: 22 -> 23
$ws.Cells.Item(7, 9).Value = 114   # I7: 113 -> 114

# Sheet 65: West Town
$ws = $wb.Worksheets.Item(65)
$ws.Cells.Item(2, 9).Value = 19   # I2: 18 -> 19
$ws.Cells.Item(4, 6).Value = 13   # F4: 12 -> 13
$ws.Cells.Item(7, 6).Value = 156   # F7: 155 -> 156
$ws.Cells.Item(7, 9).Value = 93   # I7: 92 -> 93

# Sheet 67: O'Hare
$ws = $wb.Worksheets.Item(67)
$ws.Cells.Item(2, 9).Value = 8   # I2: 7 -> 8
$ws.Cells.Item(7, 9).Value = 24   # I7: 23 -> 24

# Sheet 7: Austin
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(6, 9).Value = 257   # I6: 254 -> 257
$ws.Cells.Item(7, 9).Value = 790   # I7: 787 -> 790

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item(72)
$ws.Cells.Item(4, 9).Value = 41   # I4: 40 -> 41
$ws.Cells.Item(7, 9).Value = 78   # I7: 77 -> 78

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item(74)
$ws.Cells.Item(2, 9).Value = 53   # I2: 52 -> 53
$ws.Cells.Item(6, 9).Value = 57   # I6: 56 -> 57
$ws.Cells.Item(7, 9).Value = 164   # I7: 162 -> 164

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item(75)
$ws.Cells.Item(6, 9).Value = 51   # I6: 50 -> 51
$ws.Cells.Item(7, 9).Value = 125   # I7: 124 -> 125

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item(79)
$ws.Cells.Item(2, 9).Value = 21   # I2: 20 -> 21
$ws.Cells.Item(3, 9).Value = 19   # I3: 18 -> 19
$ws.Cells.Item(7, 9).Value = 114   # I7: 112 -> 114

# Sheet 9: Auburn Gresham
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 9).Value = 146   # I2: 145 -> 146
$ws.Cells.Item(6, 9).Value = 108   # I6: 107 -> 108
$ws.Cells.Item(7, 9).Value = 421   # I7: 419 -> 421
